$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.485.31'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '1.602.82'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.68'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.522'
$ws.Range('E6').Value = '  +7.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.88'
$ws.Range('E8').Value = '  +8.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.56'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('E10').Value = '  +1.96%  '
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0911'
$ws.Range('E12').Value = '  +2.05%  '
$ws.Range('D13').Value = '1.831.58'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').Value = '1.596.53'
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = '29.498.33'
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.536'
$ws.Range('E16').Value = '  +3.96%  '
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.43'
$ws.Range('E18').Value = '  +3.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.85'
$ws.Range('E19').Value = '  +5.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'
$ws.Range('E20').Value = '  +3.28%  '
$ws.Range('D21').Value = '0.0₃0690'
$ws.Range('E21').Value = '  +2.61%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.17'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.56'
$ws.Range('E26').Value = '  +2.29%  '
$ws.Range('E27').Value = '  +5.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.29'
$ws.Range('E28').Value = '  +3.44%  '
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0473'
$ws.Range('E31').Value = '  +2.71%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('D34').Value = '1.417.52'
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('E35').Value = '  +3.68%  '
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('E37').Value = '  +2.47%  '
$ws.Range('E38').Value = '  +5.31%  '
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.535'
$ws.Range('E41').Value = '  +3.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.96'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0479'
$ws.Range('E44').Value = '  +3.13%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '52.80'
$ws.Range('E45').Value = '  +21.53%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.792'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.58'
$ws.Range('E47').Value = '  +2.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.28'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').Value = '1.742.26'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.45'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.829'
$ws.Range('E51').Value = '  -4.55%  '
